$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new BOM row (row 5) describing the 18V-to-5V DCDC converter part.
# Set values in this particular order so shared-string entries are appended
# in the same sequence as the authoritative edit (Type, Link, Name, Package).
$ws.Range("A5").Value = "18 to 5 V DCDC"
$ws.Range("D5").Value = "https://hu.farnell.com/diodes-inc/ap1509-50sg-13/ic-buck-reg-5v-2a-8sop/dp/1825323"
$ws.Range("B5").Value = "DIODES INC. AP1509-50SG-13"
$ws.Range("C5").Value = "SOP-8L"

# Widen column B to fit the longer part name.
$ws.Columns.Item(2).ColumnWidth = 25.88671875

# Move the active selection, matching where the author left the cursor.
$ws.Range("L10").Select() | Out-Null
